{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 1 (index 0): \"Jia Ming\" -> \"Jia. Ming\" ---\nconst hits1 = body.search(\"Jia Ming\", { matchCase: true });\nhits1.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < hits1.items.length; i++) {\n  hits1.items[i].insertText(\"Jia. Ming\", \"Replace\");\n}\nawait context.sync();\n\n// --- Paragraph 2 (index 1): several small text fixes ---\nconst replacements = [\n  [\"This type | of water\", \"This. type of water\"],\n  [\"switched on. and\", \"switched.on and\"],\n  [\"passes through. the\", \"passes through the\"],\n  [\"heater, The water heater\", \"heater. The water heater\"],\n  [\"temperature of 49 \\u00b0C. :\", \"temperature of 49 \\u00b0C, :\"],\n];\nfor (const [find, repl] of replacements) {\n  const hits = body.search(find, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < hits.items.length; i++) {\n    hits.items[i].insertText(repl, \"Replace\");\n  }\n  await context.sync();\n}\n\n// --- Paragraph 3 (index 2): remove the inline picture, replace with new text ---\nconst imgParagraph = paragraphs.items[2];\nconst pics = imgParagraph.inlinePictures;\npics.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < pics.items.length; i++) {\n  pics.items[i].delete();\n}\nawait context.sync();\nimgParagraph.insertText(\n  \"(d} | Which water heater A or B is more energy saving? Explain your choice.\\u000b[11\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- Old paragraph 4 (index 3): \"(d) Which water heater...\" paragraph is removed entirely ---\nparagraphs.items[3].delete();\nawait context.sync();\n\n// --- Old paragraph 5 (now index 3 after delete): replace its text ---\nparagraphs.load(\"items\");\nawait context.sync();\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertText(\"of : i Pe\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# --- Paragraph 1: \"Jia Ming\" -> \"Jia. Ming\" ---\nReplace-Text \"Jia Ming\" \"Jia. Ming\"\n\n# --- Paragraph 2: several small text fixes ---\nReplace-Text \"This type | of water\" \"This. type of water\"\nReplace-Text \"switched on. and\" \"switched.on and\"\nReplace-Text \"passes through. the\" \"passes through the\"\nReplace-Text \"heater, The water heater\" \"heater. The water heater\"\nReplace-Text ([char]0x00B0 + \"C. :\") ([char]0x00B0 + \"C, :\")\n\n# --- Paragraph 3: the inline picture's paragraph becomes new text (picture removed) ---\n$p3 = $d.Paragraphs.Item(3)\n$p3.Range.Text = \"(d} | Which water heater A or B is more energy saving? Explain your choice.\" + [char]11 + \"[11\"\n\n# --- Paragraph 4 (old \"(d) Which water heater...\" paragraph) is removed entirely ---\n$d.Paragraphs.Item(4).Range.Delete()\n\n# --- Old paragraph 5 (now paragraph 4): replace its text ---\n$d.Paragraphs.Item(4).Range.Text = \"of : i Pe\"\n"}
